# Zeitplanung.xlsx update
# - Analyse & Design (row 15): add 0.5h on K15 (new "Besprechung Anforderungsanalyse" slot),
#   bump planned effort C15 from 1 to 1.5
# - Implementation section (rows 21-29): re-balance requirement tasks/hours,
#   drop the last two placeholder rows (310/311)
# - Update active selection / scroll position to reflect where work is happening

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Zeitplanung")

# --- Analyse & Design rows ---
$ws.Range("C15").Value = 1.5
$ws.Range("K15").Value = 0.5

# --- Implementation rows (21-29) ---
$ws.Range("B21").Value = "Anforderung NF001 & NF002"
$ws.Range("C21").Value = 12

$ws.Range("C23").Value = 10

$ws.Range("B24").Value = "Anforderung A002"
$ws.Range("C24").Value = 8

$ws.Range("B25").Value = "Anforderung A003"
$ws.Range("C25").Value = 9

$ws.Range("C26").Value = 5

$ws.Range("B27").Value = "Anforderung A005"
$ws.Range("C27").Value = 4

$ws.Range("B28").ClearContents()
$ws.Range("C28").ClearContents()

$ws.Range("B29").ClearContents()
$ws.Range("C29").ClearContents()

# --- View state: selection moved to K15, scrolled up one row ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 14
$ws.Range("K15").Select()
